$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1443736666666667
$ws.Range("H2").Value = 0.433121
$ws.Range("I2").Value = 0.7378778224885942
$ws.Range("J2").Value = 0.7378778224885942
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.038663
$ws.Range("N2").Value = 0.115989
$ws.Range("O2").Value = 0.05376113331800686
$ws.Range("P2").Value = 0.05376113331800687
$ws.Range("Q2").Value = 0.005581919074333333
$ws.Range("R2").Value = 0.050237271669
$ws.Range("S2").Value = 0.03966914798720991
$ws.Range("T2").Value = 0.03966914798720992

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1443736666666667
$ws.Range("H3").Value = 0.433121
$ws.Range("I3").Value = 0.7378778224885942
$ws.Range("J3").Value = 0.7378778224885942
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.08844233333333333
$ws.Range("N3").Value = 0.265327
$ws.Range("O3").Value = 0.1229795947880127
$ws.Range("P3").Value = 0.1229795947880127
$ws.Range("Q3").Value = 0.01276874395188889
$ws.Range("R3").Value = 0.114918695567
$ws.Range("S3").Value = 0.09074391561270849
$ws.Range("T3").Value = 0.0907439156127085

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1443736666666667
$ws.Range("H4").Value = 0.433121
$ws.Range("I4").Value = 0.7378778224885942
$ws.Range("J4").Value = 0.7378778224885942
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.074226
$ws.Range("N4").Value = 0.222678
$ws.Range("O4").Value = 0.103211698048842
$ws.Range("P4").Value = 0.103211698048842
$ws.Range("Q4").Value = 0.010716279782
$ws.Range("R4").Value = 0.096446518038
$ws.Range("S4").Value = 0.0761576230116298
$ws.Range("T4").Value = 0.07615762301162982

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1443736666666667
$ws.Range("H5").Value = 0.433121
$ws.Range("I5").Value = 0.7378778224885942
$ws.Range("J5").Value = 0.7378778224885942
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5178313333333334
$ws.Range("N5").Value = 1.553494
$ws.Range("O5").Value = 0.7200475738451385
$ws.Range("P5").Value = 0.7200475738451385
$ws.Range("Q5").Value = 0.07476120830822224
$ws.Range("R5").Value = 0.6728508747740001
$ws.Range("S5").Value = 0.531307135877046
$ws.Range("T5").Value = 0.531307135877046

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.051287
$ws.Range("H6").Value = 0.153861
$ws.Range("I6").Value = 0.2621221775114058
$ws.Range("J6").Value = 0.2621221775114058
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.038663
$ws.Range("N6").Value = 0.115989
$ws.Range("O6").Value = 0.05376113331800686
$ws.Range("P6").Value = 0.05376113331800687
$ws.Range("Q6").Value = 0.001982909281
$ws.Range("R6").Value = 0.017846183529
$ws.Range("S6").Value = 0.01409198533079694
$ws.Range("T6").Value = 0.01409198533079695

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.051287
$ws.Range("H7").Value = 0.153861
$ws.Range("I7").Value = 0.2621221775114058
$ws.Range("J7").Value = 0.2621221775114058
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.08844233333333333
$ws.Range("N7").Value = 0.265327
$ws.Range("O7").Value = 0.1229795947880127
$ws.Range("P7").Value = 0.1229795947880127
$ws.Range("Q7").Value = 0.004535941949666666
$ws.Range("R7").Value = 0.040823477547
$ws.Range("S7").Value = 0.03223567917530422
$ws.Range("T7").Value = 0.03223567917530423

$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.051287
$ws.Range("H8").Value = 0.153861
$ws.Range("I8").Value = 0.2621221775114058
$ws.Range("J8").Value = 0.2621221775114058
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.074226
$ws.Range("N8").Value = 0.222678
$ws.Range("O8").Value = 0.103211698048842
$ws.Range("P8").Value = 0.103211698048842
$ws.Range("Q8").Value = 0.003806828862
$ws.Range("R8").Value = 0.034261459758
$ws.Range("S8").Value = 0.02705407503721217
$ws.Range("T8").Value = 0.02705407503721218

$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.051287
$ws.Range("H9").Value = 0.153861
$ws.Range("I9").Value = 0.2621221775114058
$ws.Range("J9").Value = 0.2621221775114058
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5178313333333334
$ws.Range("N9").Value = 1.553494
$ws.Range("O9").Value = 0.7200475738451385
$ws.Range("P9").Value = 0.7200475738451385
$ws.Range("Q9").Value = 0.02655801559266667
$ws.Range("R9").Value = 0.239022140334
$ws.Range("S9").Value = 0.1887404379680925
$ws.Range("T9").Value = 0.1887404379680925

Write-Output "Done"
